$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("master-reg_center_device")

# Append 5 new device rows for regcntr_id 10002 (device_id 3000176-3000180)
$startRow = 157
$startDevice = 3000176

for ($i = 0; $i -lt 5; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = 10002
    $ws.Cells.Item($row, 2).Value = ($startDevice + $i)
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
}

# Match the updated view/selection state from the edit
$ws.Range("C158").Select()
